$d = $word.ActiveDocument

# Locate the paragraph that holds the "public-pixel-font" hyperlink (the
# last hyperlink paragraph in the document) by searching for its text,
# then insert a brand-new paragraph right after it.
$finder = $d.Content
$finder.Find.Execute("https://ggbot.itch.io/public-pixel-font", $true, $false,
    $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$anchorParaIndex = $finder.Paragraphs(1).Index
$anchor = $d.Paragraphs($anchorParaIndex)
$insertAt = $anchor.Range.End - 1
$breakPoint = $d.Range($insertAt, $insertAt)
$breakPoint.InsertParagraphAfter()

# Put the plain URL text into the freshly created paragraph.
$newParaIndex = $anchorParaIndex + 1
$newPara = $d.Paragraphs($newParaIndex)
$startOfNewPara = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$url = "https://zvukipro.com/"
$startOfNewPara.InsertAfter($url)

# Re-acquire the paragraph/range now that it actually contains the text,
# and turn that exact span into a hyperlink.
$newPara = $d.Paragraphs($newParaIndex)
$textRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + $url.Length)
$hyperlink = $d.Hyperlinks.Add($textRange, $url, $null, $null, $url)

# Make sure the new run uses the same character style ("Hyperlink",
# internally stored as styleId "a3") as every other link in the document.
$textRange.Style = $d.Styles.Item("a3")
